# Trade #20 closed at 2026-02-17 20:04:33 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - update aggregate stats now that trade #20 has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.8    # Current Capital
$summary.Range("B4").Value = -0.21     # Total P&L $
$summary.Range("B6").Value = 20        # Total Trades
$summary.Range("B8").Value = 10        # Losing Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 5) stats.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.8       # Capital
$status.Range("D5").Value = 20         # Trades
$status.Range("E5").Value = -0.21      # P&L $
$status.Range("F5").Value = -0.2       # P&L %
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append the new trade record (row 21) to both "All Trades" and
#    "MarketMaking" sheets - they mirror each other.
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(21, 1).Value  = 20                                     # A21 Trade #

    # Date/Time columns: force text format first so Excel's COM layer does
    # not auto-coerce the "yyyy-mm-dd" / "hh:mm:ss"-looking strings into
    # date/time serial numbers, then restore the Normal style so no stray
    # formatting is left behind on the cell.
    $ws.Cells.Item(21, 2).NumberFormat = "@"
    $ws.Cells.Item(21, 2).Value = "2026-02-17"                            # B21 Date
    $ws.Cells.Item(21, 2).Style = "Normal"

    $ws.Cells.Item(21, 3).NumberFormat = "@"
    $ws.Cells.Item(21, 3).Value = "20:04:26"                              # C21 Time
    $ws.Cells.Item(21, 3).Style = "Normal"

    $ws.Cells.Item(21, 4).Value  = "MarketMaking"                         # D21 Strategy
    $ws.Cells.Item(21, 5).Value  = "DOWN"                                 # E21 Side
    $ws.Cells.Item(21, 6).Value  = 0.03                                   # F21 Entry Price
    $ws.Cells.Item(21, 7).Value  = 0.02                                   # G21 Exit Price
    $ws.Cells.Item(21, 8).Value  = "CLOSED"                               # H21 Status
    $ws.Cells.Item(21, 9).Value  = -33.3333                               # I21 P&L %
    $ws.Cells.Item(21, 10).Value = -0.01                                  # J21 P&L $
    $ws.Cells.Item(21, 11).Value = 99.8                                   # K21 Capital After
    $ws.Cells.Item(21, 12).Value = 0                                      # L21 Entry Slippage (bps)
    $ws.Cells.Item(21, 13).Value = 0                                      # M21 Exit Slippage (bps)
    $ws.Cells.Item(21, 14).Value = 0.6                                    # N21 Confidence
    $ws.Cells.Item(21, 15).Value = "Normal spread capture: 19600 bps"     # O21 Entry Reason
    $ws.Cells.Item(21, 16).Value = "early_exit"                           # P21 Exit Reason
    $ws.Cells.Item(21, 17).Value = 0.13                                   # Q21 Duration (min)
}
